$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "44.068.72"

Set-TextValue $ws.Range("D3") "2.243.61"
Set-TextValue $ws.Range("E3") "  -1.43%  "

Set-TextValue $ws.Range("E4") "  +0.12%  "

Set-TextValue $ws.Range("D5") "315.70"
Set-TextValue $ws.Range("E5") "  -1.75%  "

Set-TextValue $ws.Range("D6") "99.38"
Set-TextValue $ws.Range("E6") "  -6.28%  "

Set-TextValue $ws.Range("D7") "0.575"
Set-TextValue $ws.Range("E7") "  -3.15%  "

Set-TextValue $ws.Range("E8") "  +0.10%  "

Set-TextValue $ws.Range("D9") "0.534"
Set-TextValue $ws.Range("E9") "  -6.62%  "

Set-TextValue $ws.Range("D10") "36.35"
Set-TextValue $ws.Range("E10") "  -6.18%  "

Set-TextValue $ws.Range("D11") "0.0823"
Set-TextValue $ws.Range("E11") "  -2.40%  "

Set-TextValue $ws.Range("D12") "7.37"
Set-TextValue $ws.Range("E12") "  -6.67%  "

Set-TextValue $ws.Range("E13") "  -2.73%  "

Set-TextValue $ws.Range("D14") "2.585.05"
Set-TextValue $ws.Range("E14") "  -1.49%  "

Set-TextValue $ws.Range("D15") "0.845"
Set-TextValue $ws.Range("E15") "  -4.48%  "

Set-TextValue $ws.Range("D16") "2.247.15"
Set-TextValue $ws.Range("E16") "  -1.52%  "

Set-TextValue $ws.Range("D17") "13.99"
Set-TextValue $ws.Range("E17") "  -4.31%  "

Set-TextValue $ws.Range("D18") "43.911.48"
Set-TextValue $ws.Range("E18") "  -0.98%  "

Set-TextValue $ws.Range("D19") "13.17"
Set-TextValue $ws.Range("E19") "  -6.32%  "

Set-TextValue $ws.Range("D20") "0.0₃0979"
Set-TextValue $ws.Range("E20") "  -2.49%  "

Set-TextValue $ws.Range("D21") "6.34"
Set-TextValue $ws.Range("E21") "  -3.08%  "

Set-TextValue $ws.Range("D22") "65.72"
Set-TextValue $ws.Range("E22") "  -1.17%  "

Set-TextValue $ws.Range("D23") "237.90"
Set-TextValue $ws.Range("E23") "  -0.54%  "

Set-TextValue $ws.Range("E24") "  -7.27%  "

Set-TextValue $ws.Range("E25") "  -8.27%  "

Set-TextValue $ws.Range("E26") "  +0.24%  "

Set-TextValue $ws.Range("D27") "10.16"
Set-TextValue $ws.Range("E27") "  -0.49%  "

Set-TextValue $ws.Range("E28") "  -4.72%  "

Set-TextValue $ws.Range("D29") "36.50"
Set-TextValue $ws.Range("E29") "  -4.87%  "

Set-TextValue $ws.Range("D30") "5.99"
Set-TextValue $ws.Range("E30") "  -8.21%  "

Set-TextValue $ws.Range("D31") "20.08"
Set-TextValue $ws.Range("E31") "  -2.72%  "

Set-TextValue $ws.Range("D32") "156.36"
Set-TextValue $ws.Range("E32") "  -4.60%  "

Set-TextValue $ws.Range("D33") "0.0841"
Set-TextValue $ws.Range("E33") "  -5.10%  "

Set-TextValue $ws.Range("B34") "LidoDAOToken"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D34") "3.30"
Set-TextValue $ws.Range("E34") "  +3.64%  "

Set-TextValue $ws.Range("B35") "WEMIXToken"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D35") "2.66"
Set-TextValue $ws.Range("E35") "  -3.33%  "

Set-TextValue $ws.Range("D36") "1.90"
Set-TextValue $ws.Range("E36") "  -6.79%  "

Set-TextValue $ws.Range("D37") "0.108"
Set-TextValue $ws.Range("E37") "  -7.24%  "

Set-TextValue $ws.Range("E38") "  -3.13%  "

Set-TextValue $ws.Range("D39") "15.42"
Set-TextValue $ws.Range("E39") "  -1.63%  "

Set-TextValue $ws.Range("D40") "3.54"
Set-TextValue $ws.Range("E40") "  -11.28%  "

Set-TextValue $ws.Range("D41") "3.99"
Set-TextValue $ws.Range("E41") "  -10.78%  "

Set-TextValue $ws.Range("D42") "0.0309"
Set-TextValue $ws.Range("E42") "  -6.09%  "

Set-TextValue $ws.Range("E43") "  +0.08%  "

Set-TextValue $ws.Range("D44") "1.704.09"
Set-TextValue $ws.Range("E44") "  -4.46%  "

Set-TextValue $ws.Range("D45") "82.75"
Set-TextValue $ws.Range("E45") "  -4.79%  "

Set-TextValue $ws.Range("E46") "  -6.32%  "

Set-TextValue $ws.Range("D47") "5.19"
Set-TextValue $ws.Range("E47") "  -5.63%  "

Set-TextValue $ws.Range("D48") "101.89"
Set-TextValue $ws.Range("E48") "  -2.42%  "

Set-TextValue $ws.Range("D49") "71.29"
Set-TextValue $ws.Range("E49") "  -4.95%  "

Set-TextValue $ws.Range("D50") "56.47"
Set-TextValue $ws.Range("E50") "  -6.42%  "

Set-TextValue $ws.Range("D51") "1.61"
Set-TextValue $ws.Range("E51") "  -5.49%  "
